$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy'
$ws.Range('G3').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat'
$ws.Range('G4').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range('G5').Value = 'Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma, Dr. Nourhan Mahmoud'
$ws.Range('G6').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G7').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range('G8').Value = 'Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda, Dr. Majorelle Magdy, Administrator'
$ws.Range('G9').Value = 'Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Gehan Adel'
$ws.Range('G10').Value = 'Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range('G11').Value = 'Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range('G13').Value = 'Dr. Mariam Nour El-Din, Dr. Safa Hany, D Wessam Atef, Dr. Omnia Mohammad, Dr. Shimaa Ashraf'
$ws.Range('G14').Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G17').Value = 'Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Esraa Mostafa, Dr. Nourhan Osama'
$ws.Range('G19').Value = 'D Mariam E. Mohammad, Dr. Sarah Mahdy'
$ws.Range('G22').Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range('G24').Value = 'Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Youstina Magdy, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Aya Emad, Dr. Ola Abd Al-Fattah'
$ws.Range('G25').Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Remon, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Aya Emad, Dr. Ola Abd Al-Fattah'
$ws.Range('G26').Value = 'Dr. Gehad Salah, Dr. Youstina Magdy'
$ws.Range('G27').Value = 'Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah'
$ws.Range('G28').Value = 'Dr. Neveen Nashaat, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range('G29').Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah'
$ws.Range('G30').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range('G31').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat'
$ws.Range('G32').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range('G33').Value = 'Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma, Dr. Nourhan Mahmoud'
$ws.Range('G34').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G35').Value = 'Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range('G36').Value = 'Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda, Dr. Majorelle Magdy, Administrator'
$ws.Range('G37').Value = 'Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Gehan Adel'
$ws.Range('G38').Value = 'Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range('G39').Value = 'Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range('G41').Value = 'Dr. Mariam Nour El-Din, Dr. Safa Hany, D Wessam Atef, Dr. Omnia Mohammad, Dr. Shimaa Ashraf'
$ws.Range('G42').Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G45').Value = 'Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Esraa Mostafa, Dr. Nourhan Osama'
$ws.Range('G47').Value = 'D Mariam E. Mohammad, Dr. Sarah Mahdy'
$ws.Range('G50').Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range('G52').Value = 'Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Monica, Dr. Youstina Magdy, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Aya Emad, Dr. Ola Abd Al-Fattah'
$ws.Range('G53').Value = 'Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Remon, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Aya Emad, Dr. Ola Abd Al-Fattah'
$ws.Range('G54').Value = 'Dr. Gehad Salah, Dr. Youstina Magdy'
$ws.Range('G55').Value = 'Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah'
$ws.Range('G56').Value = 'Dr. Neveen Nashaat, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range('G57').Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah'
